$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New monthly data rows to append after the existing last row (224)
$data = @(
    @{ Row = 225; Year = 2023; Month = 7; Day = 1; B = 622009300000; C = 0.2230848168473654; D = 138760830767.8579 },
    @{ Row = 226; Year = 2023; Month = 8; Day = 1; B = 626029700000; C = 0.2195775328268412; D = 137462057002.3275 },
    @{ Row = 227; Year = 2023; Month = 9; Day = 1; B = 640762400000; C = 0.2126709342634142; D = 136271538248.8675 }
)

# Row that holds the reference formatting (date style) to copy into column A of new rows
$templateRow = 224

foreach ($item in $data) {
    $r = $item.Row

    # Column A: date value, formatted/styled like the preceding rows
    $ws.Cells.Item($templateRow, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Cells.Item($r, 1).Value = Get-Date -Year $item.Year -Month $item.Month -Day $item.Day -Hour 0 -Minute 0 -Second 0

    # Columns B, C, D: plain numeric values, no special styling
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}

$excel.CutCopyMode = 0
